# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the per-row Fecha/Volumen/Precio values
# (columns D, M, N, O, P, S) among the existing data rows (2-20).
# Capture the current values first, then write back the permuted set so
# that reads and writes don't clobber each other mid-flight.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values to copy from source into destination)
$rowMap = @{
    2  = 19
    3  = 15
    4  = 13
    5  = 12
    6  = 10
    7  = 5
    8  = 14
    9  = 9
    10 = 16
    11 = 3
    12 = 20
    13 = 6
    14 = 8
    15 = 7
    16 = 2
    17 = 18
    18 = 11
    19 = 17
    20 = 4
}

# Columns that carry the permuted data (by column letter)
$cols = @("D", "M", "N", "O", "P", "S")

# Snapshot the existing values for every affected cell before writing anything
$snapshot = @{}
foreach ($row in 2..20) {
    foreach ($col in $cols) {
        $addr = "$col$row"
        $snapshot[$addr] = $ws.Range($addr).Value2
    }
}

# Apply the permutation using the captured snapshot values
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $srcAddr = "$col$srcRow"
        $destAddr = "$col$destRow"
        $ws.Range($destAddr).Value2 = $snapshot[$srcAddr]
    }
}
